$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: picture index 0, L/M/N/O values
$ws.Range("J1").Value = 0
$ws.Range("L1").Value = 0
$ws.Range("M1").Value = 0
$ws.Range("N1").Value = 15
$ws.Range("O1").Value = 10

# Row 2: picture index 1
$ws.Range("J2").Value = 1
$ws.Range("L2").Value = 40
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 27
$ws.Range("O2").Value = 36

# Row 3: picture index 2
$ws.Range("J3").Value = 2
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 39
$ws.Range("O3").Value = 36

# Row 4: picture index 3
$ws.Range("J4").Value = 3
$ws.Range("L4").Value = 68
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 27
$ws.Range("O4").Value = 36

# Row 5: picture index 4
$ws.Range("J5").Value = 4
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 37
$ws.Range("N5").Value = 24
$ws.Range("O5").Value = 36

# Rows 6-9 keep "PictureN" text labels, but point at Picture6..Picture9
# (their numbering is effectively unchanged: row6->Picture6, row7->Picture7, row8->Picture8, row9->Picture9)
$ws.Range("J6").Value = "Picture6"
$ws.Range("J7").Value = "Picture7"
$ws.Range("J8").Value = "Picture8"
$ws.Range("J9").Value = "Picture9"

# Update the selection to A1:I1 (matches saved sheet view: activeCell I1, sqref A1:I1)
$ws.Range("A1:I1").Select()
